$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 6: was "Vaina" -> now "APM 2.6", add D6 price, E6 formula already present (shared)
$ws.Range("A6").Value = "APM 2.6"
$ws.Range("D6").Value = 239.98

# Row 7: new row - name, qty, price
$ws.Range("A7").Value = "Helicopter"
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 335.41

# Row 8: new row - name, qty, price
$ws.Range("A8").Value = "Flatmaps"
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 164.5

# Row 9: new row - name, qty, price
$ws.Range("A9").Value = "RD900 radio"
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 203.28

# Sheet view changes: scroll so row 4 is at the top, select A10:B10
$ws.Range("A10:B10").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
